# Apply the "participant assignment rules" update:
# Add 15 new rule-table rows (24-38) enforcing a single assignee / owner / owning group
# per object type (CASE_FILE, COMPLAINT, PERSON, ORGANIZATION, DOC_REPO).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Expression text reused across several rows (note row 24 uses a variant with no space before ">1")
$exprAssignee        = "participants != null && participants.containsKey('assignee') && participants['assignee'].size() >1"
$exprAssigneeNoSpace = "participants != null && participants.containsKey('assignee') && participants['assignee'].size()>1"
$exprOwner           = "participants != null && participants.containsKey('owner') && participants['owner'].size() >1"
$exprOwningGroup     = "participants != null && participants.containsKey('owning group') && participants['owning group'].size() >1"

$msgAssignee     = "Only one Assignee allowed"
$msgOwner        = "Only one owner allowed"
$msgOwningGroup  = "Only one owning group allowed"

# Row definitions: RowNumber, RuleName(B), ObjectType(C), Expression(D), ErrorMessage(F), RowHeight
$rows = @(
    @(24, "Case File -Check if entry exists for participant Type assignee already exists",       "CASE_FILE",    $exprAssigneeNoSpace, $msgAssignee,    45),
    @(25, "Case file -Check if entry exists for participant Type owner already exists",           "CASE_FILE",    $exprOwner,           $msgOwner,       45),
    @(26, "Case file -Check if entry exists for participant Type owningGroup already exists",     "CASE_FILE",    $exprOwningGroup,     $msgOwningGroup, 60),

    @(27, "Complaints -Check if entry exists for participant Type assignee already exists",       "COMPLAINT",    $exprAssignee,        $msgAssignee,    60),
    @(28, "Complaints -Check if entry exists for participant Type owner already exists",          "COMPLAINT",    $exprOwner,           $msgOwner,       60),
    @(29, "Complaints -Check if entry exists for participant Type owningGroup already exists",    "COMPLAINT",    $exprOwningGroup,     $msgOwningGroup, 75),

    @(30, "People -Check if entry exists for participant Type assignee already exists",           "PERSON",       $exprAssignee,        $msgAssignee,    45),
    @(31, "People -Check if entry exists for participant Type owner already exists",              "PERSON",       $exprOwner,           $msgOwner,       45),
    @(32, "People -Check if entry exists for participant Type owningGroup already exists",        "PERSON",       $exprOwningGroup,     $msgOwningGroup, 60),

    @(33, "Organization -Check if entry exists for participant Type assignee already exists",     "ORGANIZATION", $exprAssignee,        $msgAssignee,    60),
    @(34, "Organization -Check if entry exists for participant Type owner already exists",        "ORGANIZATION", $exprOwner,           $msgOwner,       60),
    @(35, "Organization -Check if entry exists for participant Type owningGroup already exists",  "ORGANIZATION", $exprOwningGroup,     $msgOwningGroup, 75),

    @(36, "Documents -Check if entry exists for participant Type assignee already exists",        "DOC_REPO",     $exprAssignee,        $msgAssignee,    60),
    @(37, "Documents -Check if entry exists for participant Type owner already exists",           "DOC_REPO",     $exprOwner,           $msgOwner,       60),
    @(38, "Documents -Check if entry exists for participant Type owningGroup already exists",     "DOC_REPO",     $exprOwningGroup,     $msgOwningGroup, 75)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]   # B - Rule Name
    $ws.Cells.Item($rowNum, 3).Value = $r[2]   # C - Object Type
    $ws.Cells.Item($rowNum, 4).Value = $r[3]   # D - Expression
    $ws.Cells.Item($rowNum, 6).Value = $r[4]   # F - Error message
}

# Apply the same thin, all-around border used by the rest of the rule table (B24:F38)
$tableRange = $ws.Range("B24:F38")
$tableRange.Borders.LineStyle = 1   # xlContinuous
$tableRange.Borders.Weight = 2      # xlThin
$tableRange.Borders.ColorIndex = -4105  # xlAutomatic

# Rule-name column wraps text; the other columns keep the plain bordered style
$ws.Range("B24:B38").WrapText = $true

# Set the row heights to match the new content
foreach ($r in $rows) {
    $ws.Rows.Item($r[0]).RowHeight = $r[5]
}

# Update the view: selection and scroll position
$ws.Range("D37").Select()
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Added participant assignment rule rows 24-38"
